$wb = $excel.ActiveWorkbook

# Update the "Date" value on the Metadata sheet (B8, next to "Date" label in A8)
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-11-17T10:15:13+00:00"

# Clear the Relationship column values for rows 4 and 5 on the Mapping Table sheet
$map = $wb.Worksheets.Item("Mapping Table 0")
$map.Range("A4").Value = ""
$map.Range("A5").Value = ""
